# The workbook's active tab is the "output (2)" sheet that holds the
# query-table results (columns title / name / imageUrl / Column1.price).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header "Column1.price" to just "price". Because this cell is
# the header of a table column, editing its text also renames the bound
# table column definition. The old "Column1.price" shared string becomes
# unused and is dropped, while a new "price" shared string is created.
$ws.Range("D1").Value = "price"

# The author's last recorded selection after the edit was cell E3.
$ws.Range("E3").Select()
